# Automatic update of files.
# Bumps the "Förändrad" date (column C) from 2023-09-11 (45180) to 2023-09-12 (45181)
# for every data row in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value2
    if ($val -eq 45180) {
        $cell.Value = 45181
    }
}
